# Update "想去人数" (F column) values on several sheets to reflect refreshed
# scrape output (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value  = 2747
$ws.Range("F8").Value  = 1631
$ws.Range("F9").Value  = 7441
$ws.Range("F11").Value = 7633
$ws.Range("F14").Value = 6
$ws.Range("F15").Value = 6146
$ws.Range("F18").Value = 15
$ws.Range("F24").Value = 283
$ws.Range("F25").Value = 284
$ws.Range("F26").Value = 3620
$ws.Range("F30").Value = 257
$ws.Range("F31").Value = 1085
$ws.Range("F34").Value = 2606
$ws.Range("F35").Value = 1455
$ws.Range("F38").Value = 22
$ws.Range("F39").Value = 3250
$ws.Range("F43").Value = 898
$ws.Range("F44").Value = 479
$ws.Range("F45").Value = 1272
$ws.Range("F46").Value = 224

# --- Sheet "演出" (performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 400

# --- Sheet "全部类型" (all types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value  = 1631
$ws.Range("F12").Value = 7441
$ws.Range("F13").Value = 7633
$ws.Range("F15").Value = 6146
$ws.Range("F18").Value = 15
$ws.Range("F23").Value = 283
$ws.Range("F26").Value = 284
$ws.Range("F27").Value = 3620
$ws.Range("F32").Value = 257
$ws.Range("F35").Value = 2606
$ws.Range("F36").Value = 1455
$ws.Range("F40").Value = 3250
$ws.Range("F45").Value = 898
$ws.Range("F46").Value = 479
$ws.Range("F47").Value = 1272
$ws.Range("F48").Value = 224

$wb.Save()
